$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "76.123.85"
$ws.Range("E2").Value = "  +0.19%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.029.36"
$ws.Range("E3").Value = "  +3.90%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.19%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'197.11"
$ws.Range("E5").Value = "  -0.69%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'619.76"
$ws.Range("E6").Value = "  +4.40%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.548"
$ws.Range("E8").Value = "  -0.56%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +5.01%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.029.44"
$ws.Range("E10").Value = "  +3.49%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "'0.440"
$ws.Range("E11").Value = "  -2.01%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.47%  "

# Row 13 - Toncoin
$ws.Range("D13").Value = "'5.22"
$ws.Range("E13").Value = "  +5.85%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.580.39"
$ws.Range("E14").Value = "  +3.30%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "'28.89"
$ws.Range("E15").Value = "  +3.43%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "76.058.49"
$ws.Range("E16").Value = "  +0.05%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.0000191"
$ws.Range("E17").Value = "  +1.95%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.023.19"
$ws.Range("E18").Value = "  +3.92%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'13.46"
$ws.Range("E19").Value = "  +1.84%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'8.93"
$ws.Range("E20").Value = "  +2.48%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'379.38"
$ws.Range("E21").Value = "  +2.42%  "

# Row 22 - SuiNetwork
$ws.Range("D22").Value = "'2.35"
$ws.Range("E22").Value = "  +2.95%  "

# Row 23 - Polkadot
$ws.Range("D23").Value = "'4.36"
$ws.Range("E23").Value = "  +0.87%  "

# Row 24 - WrappedeETH
$ws.Range("D24").Value = "3.170.19"
$ws.Range("E24").Value = "  +3.67%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'72.27"
$ws.Range("E25").Value = "  -0.53%  "

# Row 26 - Dai
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.12%  "

# Row 27 - NEARProtocol
$ws.Range("D27").Value = "'4.31"
$ws.Range("E27").Value = "  +0.80%  "

# Row 28 - Aptos
$ws.Range("D28").Value = "'9.72"
$ws.Range("E28").Value = "  +0.98%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +1.54%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.37%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "'8.25"
$ws.Range("E31").Value = "  +2.18%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  +0.83%  "

# Row 33 - Bittensor
$ws.Range("D33").Value = "'491.73"
$ws.Range("E33").Value = "  -0.09%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  +4.78%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.21%  "

# Row 36 - EthereumClassic
$ws.Range("D36").Value = "'20.51"
$ws.Range("E36").Value = "  +2.00%  "

# Rows 37/38 swap: Monero <-> Kaspa (with updated data)
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.120"
$ws.Range("E37").Value = "  +9.40%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'162.16"
$ws.Range("E38").Value = "  -1.77%  "

# Row 39 - WhiteBITCoin
$ws.Range("D39").Value = "'20.01"
$ws.Range("E39").Value = "  +1.70%  "

# Row 40 - Aave
$ws.Range("D40").Value = "'190.34"
$ws.Range("E40").Value = "  +6.38%  "

# Row 41 - PolygonEcosystemToken
$ws.Range("E41").Value = "  -2.61%  "

# Row 42 - Cronos
$ws.Range("E42").Value = "  -5.79%  "

# Row 44 - RenderToken
$ws.Range("D44").Value = "'5.08"
$ws.Range("E44").Value = "  +3.22%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "'0.775"
$ws.Range("E45").Value = "  +18.14%  "

# Row 46 - OKB
$ws.Range("D46").Value = "'41.38"
$ws.Range("E46").Value = "  +3.12%  "

# Row 47 - ImmutableX
$ws.Range("E47").Value = "  +4.26%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  -0.59%  "

# Row 49 - dogwifhat
$ws.Range("D49").Value = "'2.41"
$ws.Range("E49").Value = "  +5.88%  "

# Row 50 - ARBITRUM
$ws.Range("D50").Value = "'0.595"
$ws.Range("E50").Value = "  +1.30%  "

# Row 51 - Filecoin
$ws.Range("D51").Value = "'3.86"
$ws.Range("E51").Value = "  -0.33%  "
